$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.460.21'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '3.359.63'
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '3.351.39'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  +3.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.629'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.60'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000275'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").Value = '3.882.98'
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.07%  '
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").Value = '3.356.65'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = '64.337.38'
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.986'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("E28").Value = '  +2.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '586.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.80'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("E37").Value = '  -8.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = '0.0₃0759'
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").Value = '3.096.21'
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0411'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.27%  '
